$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (prices, percentages) stay as
# literal text, matching the source data which stores these as plain text.
$cells = @("E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D23","E23","D24","E24","D25","E25","E26","D39","E39","D40","E40","D41","E41","E42","E43","D44","E44","D45","E45","D46","E46","E47","D48","E50","E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "-0.55%"
$ws.Range("D3").Value = "38.85"
$ws.Range("E3").Value = "6.82%"
$ws.Range("D4").Value = "5.103"
$ws.Range("E4").Value = "0.94%"
$ws.Range("D5").Value = "0.08079"
$ws.Range("E5").Value = "-0.48%"
$ws.Range("D6").Value = "1.949"
$ws.Range("E6").Value = "-4.19%"
$ws.Range("D7").Value = "4.192"
$ws.Range("E7").Value = "0.80%"
$ws.Range("D8").Value = "7.967"
$ws.Range("E8").Value = "1.37%"
$ws.Range("D9").Value = "0.9309"
$ws.Range("E9").Value = "0.32%"
$ws.Range("D10").Value = "0.1496"
$ws.Range("E10").Value = "4.99%"
$ws.Range("D11").Value = "0.1933"
$ws.Range("E11").Value = "0.27%"
$ws.Range("D12").Value = "0.09068"
$ws.Range("E12").Value = "-0.34%"
$ws.Range("D13").Value = "0.03497"
$ws.Range("E13").Value = "1.46%"
$ws.Range("D14").Value = "0.09776"
$ws.Range("E14").Value = "-1.35%"
$ws.Range("D15").Value = "0.001405"
$ws.Range("E15").Value = "-0.13%"
$ws.Range("D16").Value = "0.005966"
$ws.Range("E16").Value = "-4.87%"
$ws.Range("D17").Value = "3.786"
$ws.Range("E17").Value = "-1.52%"
$ws.Range("D18").Value = "3.457"
$ws.Range("E18").Value = "2.76%"
$ws.Range("D19").Value = "0.3423"
$ws.Range("E19").Value = "-0.46%"
$ws.Range("D20").Value = "0.1302"
$ws.Range("E20").Value = "0.74%"
$ws.Range("D21").Value = "4.671"
$ws.Range("E21").Value = "-2.70%"
$ws.Range("D23").Value = "0.04371"
$ws.Range("E23").Value = "-0.42%"
$ws.Range("D24").Value = "0.001237"
$ws.Range("E24").Value = "0.15%"
$ws.Range("D25").Value = "0.004284"
$ws.Range("E25").Value = "-12.90%"
$ws.Range("E26").Value = "0.08%"
$ws.Range("D39").Value = "0.02043"
$ws.Range("E39").Value = "0.76%"
$ws.Range("D40").Value = "0.05104"
$ws.Range("E40").Value = "-1.16%"
$ws.Range("D41").Value = "0.007442"
$ws.Range("E41").Value = "-0.61%"
$ws.Range("E42").Value = "1.50%"
$ws.Range("E43").Value = "-1.56%"
$ws.Range("D44").Value = "0.002122"
$ws.Range("E44").Value = "-0.39%"
$ws.Range("D45").Value = "0.009121"
$ws.Range("E45").Value = "-8.48%"
$ws.Range("D46").Value = "0.00006188"
$ws.Range("E46").Value = "-1.99%"
$ws.Range("E47").Value = "0.08%"
$ws.Range("D48").Value = "0.003103"
$ws.Range("E50").Value = "0.08%"
$ws.Range("E51").Value = "0.08%"
